$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date in column C for every data row
#    (rows 2-200) from 45182 to 45184.
for ($r = 2; $r -le 200; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# 2. Row 200 gains an explicit row height (15pt, custom height flag set).
$ws.Rows.Item(200).RowHeight = 15

# 3. Append a brand-new row 201 with the new cleaning/harvest notice record.
$ws.Cells.Item(201, 1).Value = "A 42837-2023"
$ws.Cells.Item(201, 2).Value = 45182
$ws.Cells.Item(201, 3).Value = 45184
$ws.Cells.Item(201, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(201, 5).Value = "MJÖLBY"
$ws.Cells.Item(201, 7).Value = 2.2
$ws.Cells.Item(201, 8).Value = 0
$ws.Cells.Item(201, 9).Value = 0
$ws.Cells.Item(201, 10).Value = 0
$ws.Cells.Item(201, 11).Value = 0
$ws.Cells.Item(201, 12).Value = 0
$ws.Cells.Item(201, 13).Value = 0
$ws.Cells.Item(201, 14).Value = 0
$ws.Cells.Item(201, 15).Value = 0
$ws.Cells.Item(201, 16).Value = 0
$ws.Cells.Item(201, 17).Value = 0

# Match the date-serial number format used by the rest of column B/C.
$ws.Range("B201").NumberFormat = $ws.Range("B200").NumberFormat
$ws.Range("C201").NumberFormat = $ws.Range("C200").NumberFormat

# R201 mirrors the wrap-text styled (but empty) species-name cell used on
# every other row.
$ws.Range("R201").Value = ""
$ws.Range("R201").WrapText = $true
